$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 1926.5
$ws.Range("I86").Value = 1485.0834
$ws.Range("J86").Value = 3250.75
$ws.Range("K86").Value = 1485.0834
$ws.Range("L86").Value = 3250.75
$ws.Range("M86").Value = -362.0834
$ws.Range("N86").Value = -5496.75
$ws.Range("H89").Value = 1926.5
$ws.Range("I89").Value = 1485.0834
$ws.Range("J89").Value = 3250.75
$ws.Range("K89").Value = 7425.416999999999
$ws.Range("L89").Value = 16253.75
$ws.Range("M89").Value = -1809.416999999999
$ws.Range("N89").Value = -27485.75
$ws.Range("H125").Value = 3246.375
$ws.Range("I125").Value = 2341.8
$ws.Range("J125").Value = 3657.5454
$ws.Range("K125").Value = 21076.2
$ws.Range("L125").Value = 32917.9086
$ws.Range("M125").Value = -18616.2
$ws.Range("N125").Value = -37837.9086
$ws.Range("H132").Value = 17807.73
$ws.Range("I132").Value = 20245.963
$ws.Range("J132").Value = 1044.875
$ws.Range("K132").Value = 60737.889
$ws.Range("L132").Value = 3134.625
$ws.Range("M132").Value = -58207.889
$ws.Range("N132").Value = -8194.625
$ws.Range("H141").Value = 1500
$ws.Range("I141").Value = 1500
$ws.Range("K141").Value = 4500
$ws.Range("M141").Value = 680

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5130808.5
$ws.Range("I32").Value = 6360801.5
$ws.Range("J32").Value = 21606.54
$ws.Range("K32").Value = 6360801.5
$ws.Range("L32").Value = 21606.54
$ws.Range("M32").Value = -6360514.5
$ws.Range("N32").Value = -22180.54
$ws.Range("H61").Value = 1596258.9
$ws.Range("I61").Value = 8171.64
$ws.Range("J61").Value = 4904774
$ws.Range("K61").Value = 8171.64
$ws.Range("L61").Value = 4904774
$ws.Range("M61").Value = -7959.64
$ws.Range("N61").Value = -4905198
$ws.Range("H122").Value = 1488.8334
$ws.Range("I122").Value = 869.61536
$ws.Range("J122").Value = 3098.8
$ws.Range("K122").Value = 2608.84608
$ws.Range("L122").Value = 9296.400000000001
$ws.Range("M122").Value = -158.8460800000003
$ws.Range("N122").Value = -14196.4
$ws.Range("H136").Value = 1596258.9
$ws.Range("I136").Value = 8171.64
$ws.Range("J136").Value = 4904774
$ws.Range("K136").Value = 24514.92
$ws.Range("L136").Value = 14714322
$ws.Range("M136").Value = -21964.92
$ws.Range("N136").Value = -14719422

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1767.52
$ws.Range("I20").Value = 1620.8823
$ws.Range("K20").Value = 1620.8823
$ws.Range("M20").Value = -1373.8823

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1311.2
$ws.Range("I16").Value = 889
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 889
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -602
$ws.Range("N16").Value = -3574
$ws.Range("H58").Value = 1685.069
$ws.Range("I58").Value = 1237.2
$ws.Range("J58").Value = 2164.9285
$ws.Range("K58").Value = 1237.2
$ws.Range("L58").Value = 2164.9285
$ws.Range("M58").Value = -1034.2
$ws.Range("N58").Value = -2570.9285
$ws.Range("H59").Value = 65000
$ws.Range("I59").Value = 65000
$ws.Range("K59").Value = 65000
$ws.Range("M59").Value = -63855
$ws.Range("H60").Value = 30000
$ws.Range("J60").Value = 50000
$ws.Range("L60").Value = 50000
$ws.Range("N60").Value = -51022
$ws.Range("H113").Value = 1311.2
$ws.Range("I113").Value = 889
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 889
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1281
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 6062937
$ws.Range("I132").Value = 1743.2333
$ws.Range("J132").Value = 13336369
$ws.Range("K132").Value = 5229.699900000001
$ws.Range("L132").Value = 40009107
$ws.Range("M132").Value = -2699.699900000001
$ws.Range("N132").Value = -40014167
$ws.Range("H134").Value = 2229
$ws.Range("I134").Value = 1932.963
$ws.Range("K134").Value = 5798.889
$ws.Range("M134").Value = -3263.889
$ws.Range("H136").Value = 1685.069
$ws.Range("I136").Value = 1237.2
$ws.Range("J136").Value = 2164.9285
$ws.Range("K136").Value = 3711.6
$ws.Range("L136").Value = 6494.7855
$ws.Range("M136").Value = -1161.6
$ws.Range("N136").Value = -11594.7855
$ws.Range("H141").Value = 762229.75
$ws.Range("J141").Value = 746582.3
$ws.Range("L141").Value = 746582.3
$ws.Range("N141").Value = -756942.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 136107.28
$ws.Range("I128").Value = 136107.28
$ws.Range("K128").Value = 408321.84
$ws.Range("M128").Value = -403341.84
$ws.Range("H132").Value = 2600
$ws.Range("J132").Value = 2511.2222
$ws.Range("L132").Value = 22600.9998
$ws.Range("N132").Value = -27660.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5155.6665
$ws.Range("I2").Value = 143.54546
$ws.Range("K2").Value = 143.54546
$ws.Range("M2").Value = -30.54545999999999
$ws.Range("H97").Value = 6513.4443
$ws.Range("I97").Value = 1237.9
$ws.Range("J97").Value = 13107.875
$ws.Range("K97").Value = 1237.9
$ws.Range("L97").Value = 13107.875
$ws.Range("M97").Value = -741.9000000000001
$ws.Range("N97").Value = -14099.875
$ws.Range("H122").Value = 7789.5454
$ws.Range("I122").Value = 6147.5
$ws.Range("J122").Value = 12168.333
$ws.Range("K122").Value = 18442.5
$ws.Range("L122").Value = 36504.999
$ws.Range("M122").Value = -15992.5
$ws.Range("N122").Value = -41404.999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21618.564
$ws.Range("I7").Value = 22380.95
$ws.Range("J7").Value = 20816.053
$ws.Range("K7").Value = 22380.95
$ws.Range("L7").Value = 20816.053
$ws.Range("M7").Value = -22268.95
$ws.Range("N7").Value = -21040.053
$ws.Range("H40").Value = 3099.2632
$ws.Range("I40").Value = 2147.4614
$ws.Range("K40").Value = 2147.4614
$ws.Range("M40").Value = -2011.4614
$ws.Range("H55").Value = 345
$ws.Range("I55").Value = 261
$ws.Range("K55").Value = 261
$ws.Range("M55").Value = -88
$ws.Range("H122").Value = 5627.3335
$ws.Range("I122").Value = 4829.4116
$ws.Range("J122").Value = 6341.263
$ws.Range("K122").Value = 14488.2348
$ws.Range("L122").Value = 19023.789
$ws.Range("M122").Value = -12038.2348
$ws.Range("N122").Value = -23923.789
$ws.Range("H126").Value = 21618.564
$ws.Range("I126").Value = 22380.95
$ws.Range("J126").Value = 20816.053
$ws.Range("K126").Value = 67142.85000000001
$ws.Range("L126").Value = 62448.159
$ws.Range("M126").Value = -64672.85000000001
$ws.Range("N126").Value = -67388.159
$ws.Range("H132").Value = 8960.311
$ws.Range("I132").Value = 6328.75
$ws.Range("J132").Value = 14808.223
$ws.Range("K132").Value = 18986.25
$ws.Range("L132").Value = 44424.669
$ws.Range("M132").Value = -16456.25
$ws.Range("N132").Value = -49484.669
$ws.Range("H136").Value = 4185.067
$ws.Range("I136").Value = 3709.2727
$ws.Range("J136").Value = 5493.5
$ws.Range("K136").Value = 11127.8181
$ws.Range("L136").Value = 16480.5
$ws.Range("M136").Value = -8577.8181
$ws.Range("N136").Value = -21580.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 142857650
$ws.Range("I107").Value = 588.6667
$ws.Range("K107").Value = 1766.0001
$ws.Range("M107").Value = 153.9999
$ws.Range("H136").Value = 3805.9
$ws.Range("I136").Value = 3068.5557
$ws.Range("K136").Value = 9205.667099999999
$ws.Range("M136").Value = -6655.667099999999
